$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to the new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# I0 / IF values for rows 2 through 38
$values = @(
    @(6,7),
    @(4,9),
    @(2,6),
    @(5,8),
    @(5,8),
    @(1,5),
    @(1,4),
    @(1,5),
    @(1,5),
    @(3,5),
    @(6,6),
    @(6,7),
    @(7,9),
    @(8,9),
    @(6,9),
    @(7,8),
    @(5,8),
    @(1,7),
    @(1,3),
    @(1,5),
    @(1,6),
    @(1,7),
    @(1,5),
    @(7,7),
    @(7,9),
    @(5,6),
    @(4,7),
    @(7,8),
    @(6,7),
    @(8,9),
    @(6,9),
    @(5,5),
    @(4,6),
    @(1,6),
    @(1,5),
    @(1,4),
    @(3,4)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
